$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update columns A-F for rows 2-7 (same new values across the block)
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = 114
    $ws.Cells.Item($r, 2).Value = "Anantara"
    $ws.Cells.Item($r, 3).Value = "Test1"
    $ws.Cells.Item($r, 4).Value = "Ayush "
    $ws.Cells.Item($r, 5).Value = "ayush_s@anantsol.com"
}

# Column F (provider_email) differs between rows 2-4 and 5-7
$ws.Cells.Item(2, 6).Value = "subhayan_r@anantsol.com"
$ws.Cells.Item(3, 6).Value = "subhayan_r@anantsol.com"
$ws.Cells.Item(4, 6).Value = "subhayan_r@anantsol.com"
$ws.Cells.Item(5, 6).Value = "ayush_s@anantsol.com"
$ws.Cells.Item(6, 6).Value = "ayush_s@anantsol.com"
$ws.Cells.Item(7, 6).Value = "ayush_s@anantsol.com"

# Column I (feedback_text) - unique per row
$ws.Cells.Item(2, 9).Value = "uu"
$ws.Cells.Item(3, 9).Value = "vv"
$ws.Cells.Item(4, 9).Value = "zz"
$ws.Cells.Item(5, 9).Value = "aaabb"
$ws.Cells.Item(6, 9).Value = "bbbaa"
$ws.Cells.Item(7, 9).Value = "cccdd"
